$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.519.56"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +1.04%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.874.16"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.40%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'1.010"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  +0.45%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'316.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.62%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'  +0.16%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.5081"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.48%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.3898"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -0.49%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.08372"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +1.15%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -1.01%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +0.46%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'6.217"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -0.12%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.870.86"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.42%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'20.43"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +0.75%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'7.242"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.41%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'1.010"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.27%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.00001103"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.41%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'91.21"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.18%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  +0.50%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'17.71"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.35%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  +0.15%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'5.928"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.41%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'28.553.03"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +1.03%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'11.07"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -0.14%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'2.236"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.67%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'2.087.29"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +0.52%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'161.64"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +1.03%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  +0.48%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'2.349"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -2.99%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'125.99"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +0.52%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'0.1045"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -1.47%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  +0.46%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'5.790"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -1.32%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'3.613"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -0.07%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'0.02458"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +1.02%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.06553"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +1.30%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.2161"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -0.59%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'8.857"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -2.90%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'5.085"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +2.81%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  -0.07%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'1.190"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +0.94%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.6424"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -0.40%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'11.11"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +0.21%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'1.008"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +0.17%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.6037"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +0.24%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'13.04"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.24%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'3.693"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +0.36%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'2.008"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -0.12%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'1.216"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +0.48%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'122.02"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +0.56%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'1.175"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -8.38%  "
$ws.Range("E51").Style = "Normal"
